# Insert a new "Lottery" income row above the existing data (new row 2),
# pushing salary2/Travel/Salary down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 2:4 down to 3:5 to make room for the new row.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the Lottery entry.
$ws.Range("A2").Value = "Lottery"
$ws.Range("B2").Value = 100000
$ws.Range("C2").Value = 45905.250231481485
$ws.Range("C2").NumberFormat = "m/d/yy"
